# rnaSample_2288.xlsx fix
# Commit message: "fixed harvester column in rnasamples -- holly added S.GISH to
# harvester in bioSamples"
#
# The 'harvester' column (column B) had incorrectly been filled with the same
# value used for 'rnaPreparer' ("Retrofitted_2288"). This corrects column B
# (rows 2-25, i.e. all data rows under the "harvester" header) to "S.GISH",
# without touching any of the other columns (rnaPreparer / rnaPrepMethod keep
# their original values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the harvester values for every data row.
$ws.Range("B2:B25").Value = "S.GISH"

# Re-fit column B to its new contents and leave it selected, mirroring the
# state the workbook was left in after the edit.
$ws.Columns("B").AutoFit() | Out-Null
$ws.Columns("B").Select() | Out-Null
